$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 75 - shifts existing rows 75..105 down to 76..106
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly record
$ws.Range("A75").Value = 6
$ws.Range("B75").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C75").Value = "Metropolitana"
$ws.Range("D75").Value = 44837
$ws.Range("E75").Value = 13
$ws.Range("F75").Value = 100114007
$ws.Range("G75").Value = "Jengibre"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 220
$ws.Range("K75").Value = 12000
$ws.Range("L75").Value = 14000
$ws.Range("M75").Value = 13091
$ws.Range("N75").Value = "$/caja 13 kilos"
$ws.Range("O75").Value = "Perú"
$ws.Range("P75").Value = 1007
$ws.Range("Q75").Value = 13
$ws.Range("R75").Value = "Hortaliza"
